$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 37 (old rows 37-120 shift down to 39-122)
$ws.Rows.Item(37).EntireRow.Insert()
$ws.Rows.Item(37).EntireRow.Insert()

# New row 37 data (Angeleno / Especial)
$ws.Range("A37").Value = 7
$ws.Range("B37").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C37").Value = "Ñuble"
$ws.Range("D37").Value = 45028
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 16
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103002
$ws.Range("J37").Value = "Ciruela"
$ws.Range("K37").Value = "Angeleno"
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 12000
$ws.Range("P37").Value = 12000
$ws.Range("Q37").Value = "$/bandeja 18 kilos granel"
$ws.Range("R37").Value = "Región de O'Higgins"
$ws.Range("S37").Value = 667
$ws.Range("T37").Value = 18

# New row 38 data (Angeleno / Primera)
$ws.Range("A38").Value = 7
$ws.Range("B38").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C38").Value = "Ñuble"
$ws.Range("D38").Value = 45028
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 16
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = "Frutos de hueso (carozo)"
$ws.Range("I38").Value = 100103002
$ws.Range("J38").Value = "Ciruela"
$ws.Range("K38").Value = "Angeleno"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 50
$ws.Range("N38").Value = 10000
$ws.Range("O38").Value = 10000
$ws.Range("P38").Value = 10000
$ws.Range("Q38").Value = "$/bandeja 18 kilos granel"
$ws.Range("R38").Value = "Región de O'Higgins"
$ws.Range("S38").Value = 556
$ws.Range("T38").Value = 18
